$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update revised M2 values for existing rows 419-421 ---
$ws.Range("C419:F419").Value = 5245350000000
$ws.Range("C420:F420").Value = 5322265000000
$ws.Range("C421:F421").Value = 5421638000000

# --- Append new rows 422-424 with new ECONOMICS:BRM2 data points ---
$newRows = @(
    @{ Row = 422; Date = 45108.41666666666; Value = 5501072000000 },
    @{ Row = 423; Date = 45139.41666666666; Value = 5591097000000 },
    @{ Row = 424; Date = 45170.41666666666; Value = 5656835000000 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy formatting of the row above (A/B column style carries the
    # bold/bordered/centered date style used by every data row) before
    # writing the new values so the row keeps the same look & feel.
    $ws.Range("A" + ($row - 1) + ":G" + ($row - 1)).Copy()
    $ws.Range("A" + $row + ":G" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "ECONOMICS:BRM2"
    $ws.Cells.Item($row, 3).Value = $r.Value
    $ws.Cells.Item($row, 4).Value = $r.Value
    $ws.Cells.Item($row, 5).Value = $r.Value
    $ws.Cells.Item($row, 6).Value = $r.Value
    $ws.Cells.Item($row, 7).Value = 0
}

$excel.CutCopyMode = 0
